$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Populate the brand-new shared strings in the same order the
#        original author introduced them, so the shared-string table
#        comes out in the same sequence. ---

$ws.Range("C10").Value = "pvar_a05"
$ws.Range("C11").Value = "pvar_a03"

$ws.Range("E15").Value = "Precip and streamflow - dry years 10% drier, water moved onto wet years. Hist ET."
$ws.Range("E16").Value = "Precip and streamflow - dry years 20% drier, water moved onto wet years. Hist ET."
$ws.Range("E17").Value = "Precip and streamflow - dry years 30% drier, water moved onto wet years. Hist ET."

$ws.Range("E9").Value  = "Precip - 10 major storms, no other storms. Hist ET and streams."
$ws.Range("E10").Value = "Precip - 5 major storms, no other storms. Hist ET and streams."
$ws.Range("E11").Value = "Precip - 3 major storms, no other storms. Hist ET and streams."

$ws.Range("E12").Value = "Rainy season 90% as large as hist. Hist ET and streams."
$ws.Range("E13").Value = "Rainy season 80% as large as hist. Hist ET and streams."
$ws.Range("E14").Value = "Rainy season 70% as large as hist. Hist ET and streams."

$ws.Range("D8").Value  = "SWBM, SVIHM"

# --- 2. Fill in the rest of column D (same model list every row). ---

$ws.Range("D9").Value  = "SWBM, SVIHM"
$ws.Range("D10").Value = "SWBM, SVIHM"
$ws.Range("D11").Value = "SWBM, SVIHM"
$ws.Range("D12").Value = "SWBM, SVIHM"
$ws.Range("D13").Value = "SWBM, SVIHM"
$ws.Range("D14").Value = "SWBM, SVIHM"
$ws.Range("D15").Value = "SWBM, SVIHM"
$ws.Range("D16").Value = "SWBM, SVIHM"
$ws.Range("D17").Value = "SWBM, SVIHM"

# --- 3. Re-lay the Scenario Name column (C) so the nine precip-variation
#        scenarios are grouped a05/a03 right after a10, and the rest of
#        the rows shift down to make room. ---

$ws.Range("C8").Value  = "hist"
$ws.Range("C9").Value  = "pvar_a10"
$ws.Range("C12").Value = "pvar_b90"
$ws.Range("C13").Value = "pvar_b80"
$ws.Range("C14").Value = "pvar_b70"
$ws.Range("C15").Value = "pvar_c10"
$ws.Range("C16").Value = "pvar_c20"
$ws.Range("C17").Value = "pvar_c30"

# --- 4. Column E row 8 keeps its note, just bumped to the new string slot. ---

$ws.Range("E8").Value = "Historical precip, gaps filled with ranked regression, created 2019.08.19 (leapdays now included!). Monthly ET (up2018_b). "

# --- 5. Give the new Notes entries (E9:E14) their own cell style (mirrors
#        the extra cellXfs entry that shows up in the saved workbook). ---

$ws.Range("E9:E14").Font.Bold = $true

# --- 6. Column D needs to be wider now that it holds "SWBM, SVIHM". ---

$ws.Columns("D").ColumnWidth = 15.44140625

# --- 7. Restore the selection to where the author left off. ---

$ws.Range("G25").Select()
